$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.486.45"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.659.86"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.39"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.81"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.56%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.589"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.66%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +7.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.410"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.52%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.59"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +5.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000189"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +15.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.137.11"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.249.20"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.638.45"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.81"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.92"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "360.03"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.17%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.13"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.60%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.55"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000104"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +16.14%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.27"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.69%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.23"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +8.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "542.85"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.85"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.66"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.43"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.33%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.68"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.36"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.02"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.56"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "166.76"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.74%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.51%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.34"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +8.29%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0619"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.20"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.04%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0265"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0987"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.85"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.99%  "
